$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HLT_dict")

# Add two new dictionary rows (FINISHED / FINISHED_CONT) below the existing
# key/de/en entries, for the additional "finished" pages.
# Values are entered in "de" column first, then "en" column, to mirror the
# original authoring order of the shared string table.

$ws.Range("A38").Value = "FINISHED"
$ws.Range("B38").Value = "Sie haben die Segmentierungsaufgabe beendet."
$ws.Range("A39").Value = "FINISHED_CONT"
$ws.Range("B39").Value = "Sie haben die Segmentierungsaufgabe beendet. <br> Jetzt folgen noch ein paar weitere Fragen."
$ws.Range("C38").Value = "You finished the segmentation task."
$ws.Range("C39").Value = "You finished the segmentation task. <br>  Now on to a few more questions."

# Match the formatting used by the rest of the table (vertical-top aligned
# cells, style index 2).
$ws.Range("A38:C39").VerticalAlignment = -4160

$ws.Range("B39").Select()
